$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A=3; B=44523.84065972222; C=44523.84085648148; D="IP Address"; E=100; F=16; G=$true; H=44523.84085648148; I="0LIBOL";
       J="ebola %>%`n  pivot_longer(Cases_Guinea:last_col(), names_to=`"case_death`") %>%`n  separate(name, into = c(`"case_death`", `"country`"), sep = `"_`") %>%`n  pivot_wider(names_from = case_death, values_from = value) %>%`n  drop_na" },
    @{ A=4; B=44523.8384375; C=44523.84164351852; D="IP Address"; E=100; F=276; G=$true; H=44523.84164351852; I="3hostc";
       J="ebola %>%`n  pivot_longer(``Cases_Guinea```:last_col(), names_to = 'number') %>%`n  separate(name, into = c(`"case_death`", `"country`"), sep = `"_`") %>%`n  ____(names_from = case_death, values_from = value) %>%`n  drop_na()" },
    @{ A=5; B=44523.8425; C=44523.84265046296; D="IP Address"; E=100; F=12; G=$true; H=44523.84265046296; I="0LIBOL";
       J="ebola %>%`n  pivot_longer(Cases_Guinea:last_col()) %>%`n  separate(name, into = c(`"case_death`", `"country`"), sep = `"_`") %>%`n  pivot_wider(names_from = case_death, values_from = value) %>%`n  drop_na" },
    @{ A=6; B=44523.19250000001; C=44523.20333333334; D="IP Address"; E=50; F=935; G=$false; H=44523.87898148148; I="2nesch";
       J=$null },
    @{ A=7; B=44523.83833333333; C=44523.8518287037; D="IP Address"; E=50; F=1165; G=$false; H=44523.87899305555; I="2dunic";
       J=$null },
    @{ A=8; B=44518.83646990741; C=44518.85869212963; D="Spam"; E=50; F=1920; G=$false; H=44523.87900462963; I="0garbc";
       J=$null }
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    if ($row.J -ne $null) {
        $ws.Cells.Item($r, 10).Value = $row.J
    }
    $r = $r + 1
}
